# For this workbook, within every 4-row "year group" (A/B/C/D quarter rows)
# the B-quarter and C-quarter data rows need to be swapped (A and D stay put),
# and then the F/G columns (which duplicated B/E with some rounding noise)
# need to be removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (B-row, C-row) to swap within each 4-row year block, rows 2..69.
$pairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16), @(19,20), @(23,24), @(27,28), @(31,32),
    @(35,36), @(39,40), @(43,44), @(47,48), @(51,52), @(55,56), @(59,60), @(63,64), @(67,68)
)

foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $range1 = $ws.Range("A$r1" + ":E$r1")
    $range2 = $ws.Range("A$r2" + ":E$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Drop the now-redundant F (产销率) and G (销售量) columns.
$ws.Columns("F:G").Delete()
